# "bringing corrections to journals"
# Add a new journal/time entry row (row 42) below the existing table:
#   A42 = 3/9/2024 (date, matches formatting of the date column above)
#   B42 = 2.5 (hours)
#   C42 = running total formula C41+B42 (continuing the existing pattern)
# This also causes the dependent totals in D2 (SUM(B:B)) and F2 (40*D2)
# to recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting/style from the row above so the new date cell
# (A42) matches the rest of column A (style index 1 / date number format).
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial(-4122)

# New row's data.
$ws.Range("A42").Value = 45360
$ws.Range("B42").Value = 2.5

# Build C42's running-total formula as a shared formula, the way Excel
# does when a formula column is extended down into a fresh row: enter it
# together with a temporary helper cell so the engine creates a shared
# formula group, then remove the helper row again.
$ws.Range("A43").Value = 45361
$ws.Range("B43").Value = 1
$ws.Range("C42:C43").Formula = "=C41+B42"
$ws.Rows("43:43").Delete()

# Leave the selection on the newly entered total cell, like the author did.
[void]$ws.Range("C42").Select()
